# presupuesto.xlsx - "se le agregaron contenido e informe contenido extra anadido"
# Updates a handful of existing item labels, turns the old row-25 subtotal
# formula into a plain line item, and appends a brand-new purchase block
# (rows 26-30) with its own subtotal and a grand "TOTAL ABSOLUTO" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing "logic gate" rows (8-11): drop the part-number suffix /
# reword to the generic description used going forward.
$ws.Cells.Item(8, 2).Value  = "18 COMPUERTA LOGICA AND"
$ws.Cells.Item(9, 2).Value  = "7 COMPUERTA LOGICA OR "
$ws.Cells.Item(10, 2).Value = "6 COMPUERTA LOGICA XOR "
$ws.Cells.Item(11, 2).Value = "4 COMPUERTA LOGICA NOT"

# Row 15: "OTROS" -> "OTROS VARIOS"
$ws.Cells.Item(15, 2).Value = "OTROS VARIOS"

# Row 19: relabeled to the 7-segment decoder line.
$ws.Cells.Item(19, 2).Value = "2 DECODIFICADOR 7 SEGMENTOS 74LS48"

# Row 25 used to hold the "=SUM(C22:C24)" subtotal; that subtotal moves
# down to row 28, and row 25 becomes an ordinary purchase line instead.
# Write the new numeric value *before* touching any formatting/formulas
# below it so the dependent SUM ranges pick it up correctly.
$ws.Cells.Item(25, 3).Value = 46
$ws.Cells.Item(25, 2).Value = "1 PROTOBOARD "

# --- New rows 26-30 ------------------------------------------------------
# Seed every new numeric/blank cell first...
$ws.Cells.Item(26, 3).Value = 10
$ws.Cells.Item(27, 3).Value = 24
$ws.Cells.Item(29, 3).Value = $null

# ...then stamp the row format (borders/fonts/number format) from the
# existing plain item row 23 onto each new row individually...
$ws.Range("A23:C23").Copy() | Out-Null
$ws.Range("A26:C26").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:C23").Copy() | Out-Null
$ws.Range("A27:C27").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:C23").Copy() | Out-Null
$ws.Range("A28:C28").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:C23").Copy() | Out-Null
$ws.Range("A29:C29").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:C23").Copy() | Out-Null
$ws.Range("A30:C30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($r = 26; $r -le 30; $r++) {
    $ws.Rows.Item($r).RowHeight = 17.4
}

# ...and finally fill in the text labels + formulas.
$ws.Cells.Item(26, 2).Value = "2 COMPUERTA LOGICA AND "
$ws.Cells.Item(27, 2).Value = "2 COMPUERTAS XOR"
$ws.Cells.Item(28, 2).Value = "subtotal"
$ws.Range("C28").Formula = "=SUM(C22:C27)"

$ws.Cells.Item(29, 2).Value = $null

$ws.Cells.Item(30, 2).Value = "TOTAL ABSOLUTO"
$ws.Range("C30").Formula = "=C16+C21+C28"

# --- View state: the selection now sits on the new blank spacer row.
$ws.Range("A29").Select()
